# Clear the contents of specific benchmark cells while preserving their
# existing cell formatting/style, per the target diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BENCHMARK")

$cellsToClear = @("D6", "I6", "D12", "D13", "I13", "D14")

foreach ($addr in $cellsToClear) {
    $ws.Range($addr).ClearContents()
}
